# Fixed locator_args for list. Updated simple demo.
#
# This script reproduces (as closely as the Excel COM object model allows)
# the authoring changes captured in the xml diff:
#   - tc0002 (sheet 3) becomes the active / selected sheet (activeTab=2)
#   - tr0001 (sheet 1) is no longer the selected tab
#   - tc0001 (sheet 2) selection moves from C7 to C10, and column C widens
#   - tc0002 (sheet 3) gets a new header cell K2 = "locator_args"
#   - tc0002 (sheet 3) gets a new data row 13 (q6a / btnFormat / Color),
#     selection moves to A13, and columns C & D widen to fit the new text
#   - four new shared strings are introduced: locator_args, q6a, btnFormat, Color

$wb = $excel.ActiveWorkbook

$trSheet = $wb.Worksheets.Item("tr0001")
$tcSheet1 = $wb.Worksheets.Item("tc0001")
$tcSheet2 = $wb.Worksheets.Item("tc0002")

# --- tc0001 (sheet 2): selection moves, and column C auto-widens ---
$tcSheet1.Range("C10").Select()
$tcSheet1.Columns.Item(3).ColumnWidth = 39.8

# --- tc0002 (sheet 3): new header cell in column K ---
$tcSheet2.Range("K2").Value = "locator_args"

# --- tc0002 (sheet 3): append new row 13 (mirrors the q4 "Click a button" row) ---
$tcSheet2.Range("A13").Value = "q6a"
$tcSheet2.Range("B13").Value = "Click a button"
$tcSheet2.Range("C13").Value = "selenium,click"
$tcSheet2.Range("D13").Value = "btnFormat"
$tcSheet2.Range("E13").Value = "sample,page1"
$tcSheet2.Range("K13").Value = "Color"

# Columns C & D widen to fit the (now longer) content across the sheet
$tcSheet2.Columns.Item(3).ColumnWidth = 27.8
$tcSheet2.Columns.Item(4).ColumnWidth = 27.42

# --- selection / active sheet bookkeeping ---
$tcSheet2.Range("A13").Select()
$tcSheet2.Activate()

Write-Output "applied simple_demo locator_args fix"
